# "full automation of locations and npc"
#
# A new "next_loc" column is inserted right before the existing
# "TILESIZE" column (column J), shifting TILESIZE to column K.
# The formulas in E2:H2 (which reference the TILESIZE cell) are
# automatically re-pointed at K2 by Excel's column insert behaviour.
# The new J2 cell holds the name of the castle to send the player to
# next ("castle").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J; this shifts the old column J (TILESIZE)
# to column K and updates any formulas that referenced it.
$ws.Columns("J:J").Insert()

# Header + value for the newly inserted "next_loc" column.
$ws.Range("J1").Value = "next_loc"
$ws.Range("J2").Value = "castle"

# Match the saved selection state from the diff.
$null = $ws.Range("L9").Select()
